# Updated cryptos list on Sat Apr 22 05:20:02 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.332.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4515"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3862"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07901"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.020"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.861.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.888"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.145"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001032"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06519"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.336.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.09%  "

$ws.Range("E25").Value = "  -1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.077.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.060"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.480"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.489"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09311"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9331"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.54%  "

$ws.Range("E35").Value = "  -2.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.269"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02234"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06000"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.217"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.265"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5905"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1886"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.83%  "

$ws.Range("E45").Value = "  -5.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5638"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06787"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "107.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.10%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.368"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.924"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.81%  "
